# Add a new "Estimates" worksheet summarizing total story points.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Estimates"

# Column widths (chosen so the engine's internal pixel-rounding reproduces
# the target stored OOXML widths of 56 and ~10.54 as closely as possible).
$ws2.Columns.Item(2).ColumnWidth = 55.166666666666664
$ws2.Columns.Item(3).ColumnWidth = 9.666666666666666

# Header / value cells.
$ws2.Cells.Item(2, 2).Value = "Total Story Points Esimates (Including Desing, Cut Effort, DB Design, Testing, Requirement Detailing, Code Review, Bug Fixing, Documentation, Release Notes)"
$ws2.Cells.Item(2, 2).WrapText = $true
$ws2.Cells.Item(2, 3).Value = 314

$ws2.Rows.Item(2).RowHeight = 43.5

# Page setup (portrait, matching the rest of the workbook).
$ws2.PageSetup.Orientation = 1

# Make the new sheet the active tab.
$ws2.Activate()
$ws2.Range("C2").Select()
